$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (row-major, columns B..K)
$newValues = @(
    ,@(-17.09680779032266, 1.711141381373074, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 2.265155253704675, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, 2.109722518011778, 1.999034751541709, -17.09680779032266, 3.691340983634445, -17.09680779032266, -17.09680779032266, -17.09680779032266, 1.315733113644522, -17.09680779032266)
    ,@(-17.09680779032266, 1.691454913370933, -17.09680779032266, -17.09680779032266, -17.09680779032266, 2.55757825363492, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(2.702469087481441, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, 2.09235318700485, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(3.75389564660252, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 0.4609924330557803, -17.09680779032266, 2.288812405510916)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, 2.795778011034709, -17.09680779032266, 3.430990855485744, -17.09680779032266, -17.09680779032266, -17.09680779032266, 2.504152255267535)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, 2.526189529534117, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 1.495492290426775, 1.913974342668783)
    ,@(-17.09680779032266, -17.09680779032266, 0.9430162421676996, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 1.71929149638248)
    ,@(-17.09680779032266, -17.09680779032266, 1.036399658232142, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 1.401531129635285, -17.09680779032266)
    ,@(-17.09680779032266, 2.492603798957335, 2.827820619235314, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 3.005712613803624, 2.795659429420865, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 4.321918316941216, 1.381760405209963, 2.352109046136197, -17.09680779032266)
    ,@(-17.09680779032266, -17.09680779032266, 1.443762998795143, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266, 0.9632178433582042, -17.09680779032266, -17.09680779032266)
    ,@(-17.09680779032266, 0.6988323951293964, 1.143149615393251, -17.09680779032266, 2.82407566088839, -17.09680779032266, -17.09680779032266, 0.2975962869118161, -17.09680779032266, 1.250350724694279)
    ,@(-17.09680779032266, 0.9476130121578937, -17.09680779032266, 1.599755686702485, -17.09680779032266, 1.734403565906191, -17.09680779032266, -17.09680779032266, -17.09680779032266, -17.09680779032266)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
